$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 11).Value = 5541
$ws.Cells.Item(3, 11).Value = 5672
$ws.Cells.Item(4, 11).Value = 1187
$ws.Cells.Item(5, 11).Value = 405
$ws.Cells.Item(6, 11).Value = 6322
$ws.Cells.Item(7, 11).Value = 19127

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Cells.Item(2, 11).Value = 13
$ws.Cells.Item(7, 11).Value = 42

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Cells.Item(2, 11).Value = 62
$ws.Cells.Item(7, 11).Value = 242

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(2, 11).Value = 351
$ws.Cells.Item(3, 11).Value = 385
$ws.Cells.Item(4, 11).Value = 72
$ws.Cells.Item(7, 11).Value = 1273

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Cells.Item(2, 11).Value = 144
$ws.Cells.Item(7, 11).Value = 425

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(2, 11).Value = 223
$ws.Cells.Item(6, 11).Value = 241
$ws.Cells.Item(7, 11).Value = 821

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Cells.Item(2, 11).Value = 109
$ws.Cells.Item(6, 11).Value = 75
$ws.Cells.Item(7, 11).Value = 324

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(2, 11).Value = 189
$ws.Cells.Item(3, 11).Value = 210
$ws.Cells.Item(6, 11).Value = 187
$ws.Cells.Item(7, 11).Value = 645

$ws = $wb.Worksheets.Item("New City")
$ws.Cells.Item(2, 11).Value = 142
$ws.Cells.Item(3, 11).Value = 109
$ws.Cells.Item(7, 11).Value = 438

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Cells.Item(3, 11).Value = 132
$ws.Cells.Item(7, 11).Value = 321

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Cells.Item(6, 11).Value = 28
$ws.Cells.Item(7, 11).Value = 76

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(2, 11).Value = 166
$ws.Cells.Item(5, 11).Value = 47
$ws.Cells.Item(7, 11).Value = 562
$ws.Cells.Item(8, 11).Value = 1273
$ws.Cells.Item(9, 11).Value = 80
$ws.Cells.Item(13, 11).Value = 22
$ws.Cells.Item(14, 11).Value = 100
$ws.Cells.Item(18, 11).Value = 128
$ws.Cells.Item(19, 11).Value = 559
$ws.Cells.Item(20, 11).Value = 444
$ws.Cells.Item(22, 11).Value = 53
$ws.Cells.Item(23, 11).Value = 197
$ws.Cells.Item(29, 11).Value = 1033
$ws.Cells.Item(30, 11).Value = 76
$ws.Cells.Item(31, 11).Value = 206
$ws.Cells.Item(33, 11).Value = 821
$ws.Cells.Item(34, 11).Value = 106
$ws.Cells.Item(35, 11).Value = 31
$ws.Cells.Item(36, 11).Value = 249
$ws.Cells.Item(37, 11).Value = 645
$ws.Cells.Item(40, 11).Value = 45
$ws.Cells.Item(41, 11).Value = 132
$ws.Cells.Item(42, 11).Value = 710
$ws.Cells.Item(48, 11).Value = 246
$ws.Cells.Item(51, 11).Value = 243
$ws.Cells.Item(52, 11).Value = 497
$ws.Cells.Item(53, 11).Value = 242
$ws.Cells.Item(57, 11).Value = 75
$ws.Cells.Item(63, 11).Value = 54
$ws.Cells.Item(65, 11).Value = 438
$ws.Cells.Item(66, 11).Value = 64
$ws.Cells.Item(67, 11).Value = 727
$ws.Cells.Item(68, 11).Value = 49
$ws.Cells.Item(69, 11).Value = 42
$ws.Cells.Item(76, 11).Value = 263
$ws.Cells.Item(79, 11).Value = 484
$ws.Cells.Item(80, 11).Value = 68
$ws.Cells.Item(82, 11).Value = 19
$ws.Cells.Item(83, 11).Value = 425
$ws.Cells.Item(85, 11).Value = 900
$ws.Cells.Item(86, 11).Value = 127
$ws.Cells.Item(88, 11).Value = 208
$ws.Cells.Item(89, 11).Value = 280
$ws.Cells.Item(91, 11).Value = 214
$ws.Cells.Item(94, 11).Value = 259
$ws.Cells.Item(95, 11).Value = 324
$ws.Cells.Item(96, 11).Value = 205
$ws.Cells.Item(98, 11).Value = 89
$ws.Cells.Item(99, 11).Value = 321
$ws.Cells.Item(101, 11).Value = 19127

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Cells.Item(2, 11).Value = 70
$ws.Cells.Item(7, 11).Value = 206

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(2, 11).Value = 208
$ws.Cells.Item(3, 11).Value = 260
$ws.Cells.Item(6, 11).Value = 202
$ws.Cells.Item(7, 11).Value = 727

$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(2, 11).Value = 294
$ws.Cells.Item(3, 11).Value = 370
$ws.Cells.Item(6, 11).Value = 293
$ws.Cells.Item(7, 11).Value = 1033

$ws = $wb.Worksheets.Item("Lake View")
$ws.Cells.Item(3, 11).Value = 58
$ws.Cells.Item(4, 11).Value = 35
$ws.Cells.Item(6, 11).Value = 120
$ws.Cells.Item(7, 11).Value = 246

$ws = $wb.Worksheets.Item("Chatham")
$ws.Cells.Item(2, 11).Value = 167
$ws.Cells.Item(3, 11).Value = 174
$ws.Cells.Item(6, 11).Value = 177
$ws.Cells.Item(7, 11).Value = 559

$ws = $wb.Worksheets.Item("River North")
$ws.Cells.Item(6, 11).Value = 138
$ws.Cells.Item(7, 11).Value = 263

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Cells.Item(6, 11).Value = 38
$ws.Cells.Item(7, 11).Value = 100

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Cells.Item(3, 11).Value = 27
$ws.Cells.Item(7, 11).Value = 132

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Cells.Item(2, 11).Value = 188
$ws.Cells.Item(3, 11).Value = 219
$ws.Cells.Item(6, 11).Value = 267
$ws.Cells.Item(7, 11).Value = 710

$ws = $wb.Worksheets.Item("Boystown")
$ws.Cells.Item(3, 11).Value = 8
$ws.Cells.Item(6, 11).Value = 22

$ws = $wb.Worksheets.Item("Douglas")
$ws.Cells.Item(4, 11).Value = 13
$ws.Cells.Item(7, 11).Value = 197

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Cells.Item(6, 11).Value = 88
$ws.Cells.Item(7, 11).Value = 205

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Cells.Item(2, 11).Value = 55
$ws.Cells.Item(3, 11).Value = 100
$ws.Cells.Item(7, 11).Value = 214

$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(3, 11).Value = 154
$ws.Cells.Item(7, 11).Value = 484

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Cells.Item(3, 11).Value = 144
$ws.Cells.Item(7, 11).Value = 444

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Cells.Item(5, 11).Value = 4
$ws.Cells.Item(7, 11).Value = 128

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Cells.Item(2, 11).Value = 99
$ws.Cells.Item(7, 11).Value = 249

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Cells.Item(2, 11).Value = 190
$ws.Cells.Item(6, 11).Value = 151
$ws.Cells.Item(7, 11).Value = 562

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Cells.Item(2, 11).Value = 38
$ws.Cells.Item(7, 11).Value = 106

$ws = $wb.Worksheets.Item("West Loop")
$ws.Cells.Item(2, 11).Value = 71
$ws.Cells.Item(7, 11).Value = 259

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Cells.Item(2, 11).Value = 15
$ws.Cells.Item(7, 11).Value = 89

$ws = $wb.Worksheets.Item("North Center")
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(7, 11).Value = 64

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Cells.Item(4, 11).Value = 4
$ws.Cells.Item(7, 11).Value = 31

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Cells.Item(6, 11).Value = 18
$ws.Cells.Item(7, 11).Value = 80

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Cells.Item(4, 11).Value = 16
$ws.Cells.Item(6, 11).Value = 53
$ws.Cells.Item(7, 11).Value = 166

$ws = $wb.Worksheets.Item("United Center")
$ws.Cells.Item(6, 11).Value = 85
$ws.Cells.Item(7, 11).Value = 208

$ws = $wb.Worksheets.Item("Uptown")
$ws.Cells.Item(2, 11).Value = 78
$ws.Cells.Item(3, 11).Value = 87
$ws.Cells.Item(6, 11).Value = 83
$ws.Cells.Item(7, 11).Value = 280

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Cells.Item(6, 11).Value = 22
$ws.Cells.Item(7, 11).Value = 47

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Cells.Item(6, 11).Value = 32
$ws.Cells.Item(7, 11).Value = 127

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Cells.Item(3, 11).Value = 65
$ws.Cells.Item(4, 11).Value = 26
$ws.Cells.Item(6, 11).Value = 80
$ws.Cells.Item(7, 11).Value = 243

$ws = $wb.Worksheets.Item("North Park")
$ws.Cells.Item(2, 11).Value = 21
$ws.Cells.Item(7, 11).Value = 49

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Cells.Item(2, 11).Value = 21
$ws.Cells.Item(7, 11).Value = 75

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Cells.Item(4, 11).Value = 21
$ws.Cells.Item(6, 11).Value = 67

$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(2, 11).Value = 294
$ws.Cells.Item(3, 11).Value = 306
$ws.Cells.Item(6, 11).Value = 223
$ws.Cells.Item(7, 11).Value = 900

$ws = $wb.Worksheets.Item("Clearing")
$ws.Cells.Item(6, 11).Value = 9
$ws.Cells.Item(7, 11).Value = 53

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(6, 11).Value = 19

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Cells.Item(6, 11).Value = 33
$ws.Cells.Item(7, 11).Value = 68

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Cells.Item(2, 11).Value = 16
$ws.Cells.Item(7, 11).Value = 45

$ws = $wb.Worksheets.Item("Little Village")
$ws.Cells.Item(2, 11).Value = 134
$ws.Cells.Item(3, 11).Value = 138
$ws.Cells.Item(7, 11).Value = 497
